$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the note text in E5 (shared string used by E5) to add extra detail
$ws.Range("E5").Value = "Upload to Code pen at this point; extra time spent debugging on browsers…"

# Row 5 grows taller (wrapped note text now needs two lines)
$ws.Rows.Item(5).RowHeight = 29

# Time spent values change for rows 5 and 6, and a new value appears in row 7
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 4
$ws.Range("D7").Value = 0.5

# Highlight the "Determine if needed at all" feature row with an accent fill
$ws.Range("A8").Interior.ThemeColor = 6

# Update the selected cell shown when the workbook is next opened
$null = $ws.Range("A11").Select()
